# Update the phone number ("Numero") for David João (row 3) on the
# contacts sheet, and move the active cell selection as recorded in the
# saved workbook view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# David João's contact number
$ws.Range("B3").Value = 919659339

# Reflect the cursor position that was active when the workbook was saved
$ws.Range("C14").Select() | Out-Null
